$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.910.78"
$ws.Range("E2").Value = "  -0.24%  "

Set-TextValue $ws.Range("D3") "3.520.49"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "588.09"
$ws.Range("E5").Value = "  -1.69%  "

Set-TextValue $ws.Range("D6") "134.16"
$ws.Range("E6").Value = "  -0.26%  "

Set-TextValue $ws.Range("D7") "3.518.95"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("E8").Value = "  +0.10%  "

Set-TextValue $ws.Range("D9") "0.490"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("E10").Value = "  +1.92%  "

Set-TextValue $ws.Range("D11") "7.15"
$ws.Range("E11").Value = "  +2.56%  "

$ws.Range("E12").Value = "  -0.11%  "

Set-TextValue $ws.Range("D13") "4.125.82"
$ws.Range("E13").Value = "  -0.42%  "

Set-TextValue $ws.Range("D14") "27.71"
$ws.Range("E14").Value = "  +2.89%  "

$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("E16").Value = "  +0.53%  "

Set-TextValue $ws.Range("D17") "3.525.60"
$ws.Range("E17").Value = "  -0.60%  "

Set-TextValue $ws.Range("D18") "64.892.26"
$ws.Range("E18").Value = "  +0.62%  "

Set-TextValue $ws.Range("D19") "10.03"
$ws.Range("E19").Value = "  +1.04%  "

Set-TextValue $ws.Range("D20") "14.29"
$ws.Range("E20").Value = "  -0.53%  "

Set-TextValue $ws.Range("D21") "5.69"
$ws.Range("E21").Value = "  -2.45%  "

Set-TextValue $ws.Range("D22") "391.12"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("E23").Value = "  -0.03%  "

Set-TextValue $ws.Range("D24") "3.667.49"
$ws.Range("E24").Value = "  -0.48%  "

Set-TextValue $ws.Range("D25") "74.35"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("E27").Value = "  -1.58%  "

Set-TextValue $ws.Range("D28") "1.58"
$ws.Range("E28").Value = "  +2.75%  "

Set-TextValue $ws.Range("D29") "7.49"
$ws.Range("E29").Value = "  -2.77%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  -0.08%  "

Set-TextValue $ws.Range("D31") "2.28"
$ws.Range("E31").Value = "  -0.30%  "

Set-TextValue $ws.Range("D32") "8.28"
$ws.Range("E32").Value = "  -2.37%  "

Set-TextValue $ws.Range("D33") "3.528.73"
$ws.Range("E33").Value = "  -0.40%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D34") "1.00"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D35") "23.98"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("E37").Value = "  +3.18%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D38") "171.53"
$ws.Range("E38").Value = "  +1.61%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D39") "5.20"
$ws.Range("E39").Value = "  +4.04%  "

Set-TextValue $ws.Range("D40") "6.96"
$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("E41").Value = "  +0.83%  "

Set-TextValue $ws.Range("D42") "0.819"
$ws.Range("E42").Value = "  -0.62%  "

Set-TextValue $ws.Range("D43") "26.54"
$ws.Range("E43").Value = "  +1.70%  "

Set-TextValue $ws.Range("D44") "1.24"
$ws.Range("E44").Value = "  +0.35%  "

$ws.Range("E45").Value = "  +0.05%  "

Set-TextValue $ws.Range("D46") "42.24"
$ws.Range("E46").Value = "  -1.69%  "

Set-TextValue $ws.Range("D47") "4.42"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("E48").Value = "  +0.72%  "

Set-TextValue $ws.Range("D49") "2.500.35"
$ws.Range("E49").Value = "  +1.55%  "

Set-TextValue $ws.Range("D50") "6.87"
$ws.Range("E50").Value = "  -0.65%  "

Set-TextValue $ws.Range("D51") "0.902"
$ws.Range("E51").Value = "  +3.02%  "
